$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns B:F (rows 2-25)
$dataBF = New-Object 'object[,]' 24,5
$dataBF[0,0] = 1.02
$dataBF[0,1] = 1.106290348459799
$dataBF[0,2] = 1.10526812937073
$dataBF[0,3] = 1.118335117395546
$dataBF[0,4] = 1.121417665009357
$dataBF[1,0] = 1.02
$dataBF[1,1] = 1.107778547123415
$dataBF[1,2] = 1.10650126255026
$dataBF[1,3] = 1.119766051408021
$dataBF[1,4] = 1.122811299200462
$dataBF[2,0] = 1.02
$dataBF[2,1] = 1.108740496940739
$dataBF[2,2] = 1.107298135276957
$dataBF[2,3] = 1.120691221117783
$dataBF[2,4] = 1.123712293816658
$dataBF[3,0] = 1.02
$dataBF[3,1] = 1.109144663755097
$dataBF[3,2] = 1.107632894230493
$dataBF[3,3] = 1.121079990326009
$dataBF[3,4] = 1.124090890115649
$dataBF[4,0] = 1.02
$dataBF[4,1] = 1.109212511422045
$dataBF[4,2] = 1.107689087410577
$dataBF[4,3] = 1.12114525643998
$dataBF[4,4] = 1.124154447578319
$dataBF[5,0] = 1.02
$dataBF[5,1] = 1.108745898360943
$dataBF[5,2] = 1.107302609305101
$dataBF[5,3] = 1.12069641654103
$dataBF[5,4] = 1.123717353347377
$dataBF[6,0] = 1.02
$dataBF[6,1] = 1.106793505738466
$dataBF[6,2] = 1.105685091855986
$dataBF[6,3] = 1.118818864592032
$dataBF[6,4] = 1.121888814629009
$dataBF[7,0] = 1.02
$dataBF[7,1] = 1.103345107078862
$dataBF[7,2] = 1.102826608777834
$dataBF[7,3] = 1.115504480117725
$dataBF[7,4] = 1.118660503980801
$dataBF[8,0] = 1.02
$dataBF[8,1] = 1.101040400424339
$dataBF[8,2] = 1.100915161569562
$dataBF[8,3] = 1.113290603651218
$dataBF[8,4] = 1.116503826009914
$dataBF[9,0] = 1.02
$dataBF[9,1] = 1.100040983937518
$dataBF[9,2] = 1.1000860490634
$dataBF[9,3] = 1.112330881831303
$dataBF[9,4] = 1.115568831088668
$dataBF[10,0] = 1.02
$dataBF[10,1] = 1.099669529040356
$dataBF[10,2] = 1.099777857586486
$dataBF[10,3] = 1.111974227053935
$dataBF[10,4] = 1.115221355163219
$dataBF[11,0] = 1.02
$dataBF[11,1] = 1.099749217771351
$dataBF[11,2] = 1.099843975844547
$dataBF[11,3] = 1.11205073857033
$dataBF[11,4] = 1.11529589804395
$dataBF[12,0] = 1.019999999999999
$dataBF[12,1] = 1.100010284035121
$dataBF[12,2] = 1.100060578423124
$dataBF[12,3] = 1.112301404168417
$dataBF[12,4] = 1.115540112269779
$dataBF[13,0] = 1.02
$dataBF[13,1] = 1.100171105263361
$dataBF[13,2] = 1.100194004819698
$dataBF[13,3] = 1.11245582462037
$dataBF[13,4] = 1.115690557054891
$dataBF[14,0] = 1.02
$dataBF[14,1] = 1.101106696955397
$dataBF[14,2] = 1.100970156172973
$dataBF[14,3] = 1.113354273542866
$dataBF[14,4] = 1.116565854028677
$dataBF[15,0] = 1.02
$dataBF[15,1] = 1.10169317201879
$dataBF[15,2] = 1.101456625249166
$dataBF[15,3] = 1.113917548229309
$dataBF[15,4] = 1.117114595465161
$dataBF[16,0] = 1.02
$dataBF[16,1] = 1.102035112307957
$dataBF[16,2] = 1.10174023554152
$dataBF[16,3] = 1.114245991243304
$dataBF[16,4] = 1.11743455756244
$dataBF[17,0] = 1.02
$dataBF[17,1] = 1.102151681465091
$dataBF[17,2] = 1.101836915924836
$dataBF[17,3] = 1.114357964104521
$dataBF[17,4] = 1.117543637977995
$dataBF[18,0] = 1.02
$dataBF[18,1] = 1.101630263400174
$dataBF[18,2] = 1.101404446126014
$dataBF[18,3] = 1.11385712515269
$dataBF[18,4] = 1.117055732072123
$dataBF[19,0] = 1.02
$dataBF[19,1] = 1.099933412855428
$dataBF[19,2] = 1.09999680053863
$dataBF[19,3] = 1.112227594178617
$dataBF[19,4] = 1.115468202206902
$dataBF[20,0] = 1.02
$dataBF[20,1] = 1.098865217921957
$dataBF[20,2] = 1.099110469776692
$dataBF[20,3] = 1.111202048787728
$dataBF[20,4] = 1.114469031038889
$dataBF[21,0] = 1.02
$dataBF[21,1] = 1.099431615602275
$dataBF[21,2] = 1.099580454373491
$dataBF[21,3] = 1.111745806189552
$dataBF[21,4] = 1.114998810060215
$dataBF[22,0] = 1.02
$dataBF[22,1] = 1.101658689532133
$dataBF[22,2] = 1.101428024055247
$dataBF[22,3] = 1.113884428066743
$dataBF[22,4] = 1.117082330244408
$dataBF[23,0] = 1.02
$dataBF[23,1] = 1.104237591077428
$dataBF[23,2] = 1.103566596026985
$dataBF[23,3] = 1.116362059871168
$dataBF[23,4] = 1.119495865016724
$ws.Range("B2:F25").Value = $dataBF

# Update columns I:N (rows 2-25)
$dataIN = New-Object 'object[,]' 24,6
$dataIN[0,0] = 1.074247237397167
$dataIN[0,1] = 1.111045055476491
$dataIN[0,2] = 1.107878231033273
$dataIN[0,3] = 1.120913104722987
$dataIN[0,4] = 1.123988193723218
$dataIN[0,5] = 1.112622866893917
$dataIN[1,0] = 1.074793307814864
$dataIN[1,1] = 1.112201630391769
$dataIN[1,2] = 1.10893320991316
$dataIN[1,3] = 1.122167762833439
$dataIN[1,4] = 1.125206177156017
$dataIN[1,5] = 1.113781084278235
$dataIN[2,0] = 1.0751447419557
$dataIN[2,1] = 1.112948477641468
$dataIN[2,2] = 1.109614229123079
$dataIN[2,3] = 1.122978327981258
$dataIN[2,4] = 1.125992964019073
$dataIN[2,5] = 1.114528992136693
$dataIN[3,0] = 1.075292029978086
$dataIN[3,1] = 1.113262089474045
$dataIN[3,2] = 1.109900145190138
$dataIN[3,3] = 1.123318787233256
$dataIN[3,4] = 1.126323415378879
$dataIN[3,5] = 1.114843049334044
$dataIN[4,0] = 1.075316733678327
$dataIN[4,1] = 1.113314725139408
$dataIN[4,2] = 1.109948129331715
$dataIN[4,3] = 1.123375934247115
$dataIN[4,4] = 1.12637888134412
$dataIN[4,5] = 1.11489575974809
$dataIN[5,0] = 1.075146711809954
$dataIN[5,1] = 1.112952669557894
$dataIN[5,2] = 1.109618051052983
$dataIN[5,3] = 1.122982878396573
$dataIN[5,4] = 1.12599738075385
$dataIN[5,5] = 1.114533190006121
$dataIN[6,0] = 1.074432181334188
$dataIN[6,1] = 1.111436245738105
$dataIN[6,2] = 1.10823510451914
$dataIN[6,3] = 1.121337390704495
$dataIN[6,4] = 1.124400095018101
$dataIN[6,5] = 1.113014612690572
$dataIN[7,0] = 1.073158353872771
$dataIN[7,1] = 1.108752161855624
$dataIN[7,2] = 1.105785565885104
$dataIN[7,3] = 1.11842779697785
$dataIN[7,4] = 1.121575084686998
$dataIN[7,5] = 1.11032671710111
$dataIN[8,0] = 1.072299080696747
$dataIN[8,1] = 1.106954452556401
$dataIN[8,2] = 1.104143804936414
$dataIN[8,3] = 1.116481023979158
$dataIN[8,4] = 1.119684473046075
$dataIN[8,5] = 1.108526454848482
$dataIN[9,0] = 1.071924587528597
$dataIN[9,1] = 1.106173987304502
$dataIN[9,2] = 1.103430774751479
$dataIN[9,3] = 1.115636312318235
$dataIN[9,4] = 1.118864028301924
$dataIN[9,5] = 1.107744881246405
$dataIN[10,0] = 1.071785117227024
$dataIN[10,1] = 1.105883774894668
$dataIN[10,2] = 1.103165597498438
$dataIN[10,3] = 1.115322280627729
$dataIN[10,4] = 1.118559002920705
$dataIN[10,5] = 1.10745425670166
$dataIN[11,0] = 1.071815050729169
$dataIN[11,1] = 1.105946040689597
$dataIN[11,2] = 1.103222493757349
$dataIN[11,3] = 1.115389653720633
$dataIN[11,4] = 1.118624444465944
$dataIN[11,5] = 1.10751661092116
$dataIN[12,0] = 1.071913066371036
$dataIN[12,1] = 1.10615000467096
$dataIN[12,2] = 1.103408861798601
$dataIN[12,3] = 1.115610359865178
$dataIN[12,4] = 1.118838820465085
$dataIN[12,5] = 1.107720864554772
$dataIN[13,0] = 1.071973408341127
$dataIN[13,1] = 1.10627563198618
$dataIN[13,2] = 1.103523645937283
$dataIN[13,3] = 1.115746308494907
$dataIN[13,4] = 1.118970867895138
$dataIN[13,5] = 1.107846670275189
$dataIN[14,0] = 1.072323883289911
$dataIN[14,1] = 1.107006205897265
$dataIN[14,2] = 1.104191080949454
$dataIN[14,3] = 1.116537047424296
$dataIN[14,4] = 1.119738884912175
$dataIN[14,5] = 1.108578281685026
$dataIN[15,0] = 1.072543076389851
$dataIN[15,1] = 1.107463924100499
$dataIN[15,2] = 1.104609169640851
$dataIN[15,3] = 1.117032585471619
$dataIN[15,4] = 1.120220156509317
$dataIN[15,5] = 1.10903664990061
$dataIN[16,0] = 1.07267069450907
$dataIN[16,1] = 1.107730706426117
$dataIN[16,2] = 1.104852827868739
$dataIN[16,3] = 1.117321456353633
$dataIN[16,4] = 1.120500701046511
$dataIN[16,5] = 1.10930381108773
$dataIN[17,0] = 1.072714169471788
$dataIN[17,1] = 1.10782163906556
$dataIN[17,2] = 1.104935874267218
$dataIN[17,3] = 1.117419925526217
$dataIN[17,4] = 1.12059633030217
$dataIN[17,5] = 1.109394872861951
$dataIN[18,0] = 1.072519583218911
$dataIN[18,1] = 1.107414835682022
$dataIN[18,2] = 1.104564334007959
$dataIN[18,3] = 1.116979436356286
$dataIN[18,4] = 1.120168538561211
$dataIN[18,5] = 1.108987491770947
$dataIN[19,0] = 1.07188421335855
$dataIN[19,1] = 1.106089951043759
$dataIN[19,2] = 1.103353990098112
$dataIN[19,3] = 1.115545374856377
$dataIN[19,4] = 1.118775699731818
$dataIN[19,5] = 1.107660725644531
$dataIN[20,0] = 1.071482607068619
$dataIN[20,1] = 1.105255130555175
$dataIN[20,2] = 1.102591109634491
$dataIN[20,3] = 1.11464216901573
$dataIN[20,4] = 1.117898368906888
$dataIN[20,5] = 1.106824719615099
$dataIN[21,0] = 1.071695708486389
$dataIN[21,1] = 1.105697858409706
$dataIN[21,2] = 1.102995707772581
$dataIN[21,3] = 1.115121124822358
$dataIN[21,4] = 1.118363611914554
$dataIN[21,5] = 1.107268076193964
$dataIN[22,0] = 1.072530199493301
$dataIN[22,1] = 1.107437017234238
$dataIN[22,2] = 1.104584593936834
$dataIN[22,3] = 1.117003452673406
$dataIN[22,4] = 1.12019186302511
$dataIN[22,5] = 1.109009704823513
$dataIN[23,0] = 1.073489430079934
$dataIN[23,1] = 1.109447507383472
$dataIN[23,2] = 1.10642035106041
$dataIN[23,3] = 1.119181216323551
$dataIN[23,4] = 1.122306678106413
$dataIN[23,5] = 1.111023050099364
$ws.Range("I2:N25").Value = $dataIN

Write-Output "Updated vm_pu values for 380 kV case"
